$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary fields ---
$ws.Range("E11").Value = 803790      # Valor Mora total
$ws.Range("C13").Value = 6           # Cant. Trabajadores
$ws.Range("F13").Value = 14          # Cant. Periodos

# --- Detail table (rows 16-33) ---
# Row 16
$ws.Range("C16").Value = "1052082914"
$ws.Range("D16").Value = "DAVID ORLANDO TORRES ALVAREZ"
$ws.Range("E16").Value = "1610"
$ws.Range("F16").Value = 27580
$ws.Range("G16").Value = 828116

# Row 17
$ws.Range("C17").Value = "1052082914"
$ws.Range("D17").Value = "DAVID ORLANDO TORRES ALVAREZ"
$ws.Range("E17").Value = "1611"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 737717

# Row 18
$ws.Range("C18").Value = "1052082914"
$ws.Range("D18").Value = "DAVID ORLANDO TORRES ALVAREZ"
$ws.Range("E18").Value = "1611"
$ws.Range("F18").Value = 27580
$ws.Range("G18").Value = 828116

# Row 19
$ws.Range("C19").Value = "1052082914"
$ws.Range("D19").Value = "DAVID ORLANDO TORRES ALVAREZ"
$ws.Range("E19").Value = "1612"
$ws.Range("F19").Value = 27580
$ws.Range("G19").Value = 828116

# Row 20
$ws.Range("C20").Value = "73149156"
$ws.Range("D20").Value = "ALEJANDRO ANTONIO JAVE ROMERO"
$ws.Range("E20").Value = "1701"
$ws.Range("F20").Value = 27578
$ws.Range("G20").Value = 689455

# Row 21
$ws.Range("C21").Value = "73149156"
$ws.Range("D21").Value = "ALEJANDRO ANTONIO JAVE ROMERO"
$ws.Range("E21").Value = "1701"
$ws.Range("F21").Value = 27578
$ws.Range("G21").Value = 737717

# Row 22
$ws.Range("C22").Value = "73149156"
$ws.Range("D22").Value = "ALEJANDRO ANTONIO JAVE ROMERO"
$ws.Range("E22").Value = "1701"
$ws.Range("F22").Value = 27578
$ws.Range("G22").Value = 877803

# Row 23
$ws.Range("C23").Value = "8853283"
$ws.Range("D23").Value = "MANUEL OSVALDO BARRAGAN ESCOBAR"
$ws.Range("E23").Value = "1701"
$ws.Range("F23").Value = 27578
$ws.Range("G23").Value = 781242

# Row 24
$ws.Range("C24").Value = "1052082914"
$ws.Range("D24").Value = "DAVID ORLANDO TORRES ALVAREZ"
$ws.Range("E24").Value = "1702"
$ws.Range("F24").Value = 27580
$ws.Range("G24").Value = 828116

# Row 25
$ws.Range("C25").Value = "1052082914"
$ws.Range("D25").Value = "DAVID ORLANDO TORRES ALVAREZ"
$ws.Range("E25").Value = "1704"
$ws.Range("F25").Value = 27580
$ws.Range("G25").Value = 828116

# Row 33 (rows 26-32 are unchanged - still RAFAEL ANTONIO CARRASQUILLA HERNANDEZ)
$ws.Range("C33").Value = "1002250377"
$ws.Range("D33").Value = "RAFAEL ANTONIO CARRASQUILLA HERNANDEZ"
$ws.Range("E33").Value = "2508"
$ws.Range("F33").Value = 66000
$ws.Range("G33").Value = 1650000
